$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.24168
$ws.Range("H2").Value = 0.72504
$ws.Range("M2").Value = 1.918906333333333
$ws.Range("N2").Value = 5.756718999999999
$ws.Range("O2").Value = 0.006524019162508824
$ws.Range("P2").Value = 0.006524019162508824
$ws.Range("Q2").Value = 0.46376128264
$ws.Range("R2").Value = 4.17385154376
$ws.Range("S2").Value = 0.006524019162508824
$ws.Range("T2").Value = 0.006524019162508824

# Row 3
$ws.Range("G3").Value = 0.24168
$ws.Range("H3").Value = 0.72504
$ws.Range("O3").Value = 0.6163557430885885
$ws.Range("P3").Value = 0.6163557430885885
$ws.Range("Q3").Value = 43.81377841744001
$ws.Range("R3").Value = 394.32400575696
$ws.Range("S3").Value = 0.6163557430885885
$ws.Range("T3").Value = 0.6163557430885885

# Row 4
$ws.Range("G4").Value = 0.24168
$ws.Range("H4").Value = 0.72504
$ws.Range("M4").Value = 29.04767233333333
$ws.Range("N4").Value = 87.143017
$ws.Range("O4").Value = 0.09875811426384234
$ws.Range("P4").Value = 0.09875811426384236
$ws.Range("Q4").Value = 7.02024144952
$ws.Range("R4").Value = 63.18217304568
$ws.Range("S4").Value = 0.09875811426384234
$ws.Range("T4").Value = 0.09875811426384236

# Row 5
$ws.Range("G5").Value = 0.24168
$ws.Range("H5").Value = 0.72504
$ws.Range("M5").Value = 81.87450533333333
$ws.Range("N5").Value = 245.623516
$ws.Range("O5").Value = 0.2783621234850603
$ws.Range("P5").Value = 0.2783621234850603
$ws.Range("Q5").Value = 19.78743044896
$ws.Range("R5").Value = 178.08687404064
$ws.Range("S5").Value = 0.2783621234850603
$ws.Range("T5").Value = 0.2783621234850603
